# ENTERCAP resumo tributario - adds the "Tributos sobre vendas" lookup sheet
# and wires the sales-tax G19 calc to it; zeroes out the old hard-coded
# "Custo Fixo" percentage (H22) now that it's computed via VLOOKUP.

$wb = $excel.ActiveWorkbook

# --- 1. New sheet "Tributos sobre vendas", placed after the last tab -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTrib = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsTrib.Name = "Tributos sobre vendas"

$wsTrib.Range("A1").Value = "FOLHA (Total Período)"
$wsTrib.Range("A2").Value = "Receita com ST"

# matches the fitted content width used on the "Dados" tab
$wsTrib.Columns.Item(1).ColumnWidth = 30

$wsTrib.Range("A21").Select() | Out-Null

# --- 2. "Apresentação" sheet: custo fixo agora calculado ------------------
$wsA = $wb.Worksheets.Item("Apresentação")

# H22 was a hard-coded 40% ("Custo Fixo - Teórico" override); now starts at 0
# since F22 is derived via VLOOKUP against the FATURAMENTO row.
$wsA.Range("H22").Value = 0

# F22 used to be H22*F14; now looked up from the FATURAMENTO* row of the table.
$wsA.Range("F22").Formula = "=VLOOKUP(""FATURAMENTO*"",B13:F50,5) * H22"

# G19 ("% TRIBUTOS SOBRE VENDAS" folha bucket) now sums matching rows from
# the new "Tributos sobre vendas" lookup sheet instead of the empty G20:G20.
$wsA.Range("G19").FormulaArray = "=SUMPRODUCT(SUMIF( B13:B18, 'Tributos sobre vendas'!A1:A40, G13:G18))"

$wsA.Range("J19").Select() | Out-Null
